# Burn up/down chart update - oppdatert t.o.m. 7/3/14
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Update the "Dag 8" (L column) actuals ---------------------------------
$ws.Range("L3").Value = 2
$ws.Range("L4").Value = 2
$ws.Range("L6").Value = 2.5
$ws.Range("L16").Value = 5
$ws.Range("L18").Value = 5.5
$ws.Range("L21").Value = 13

# Re-enter the totals row formula across E23:L23 in one shot so Excel turns
# it into a shared formula group (matches fill-right / re-enter behaviour).
$ws.Range("E23:L23").Formula = "=SUM(E1:E21)"

$wb.Application.Calculate()

# --- Update the burn-up chart ------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valAx = $chart.Axes(2)
$valAx.MaximumScale = 200

# --- Move the active selection on the sheet ---------------------------------
$ws.Range("M18").Select()
